$d = $word.ActiveDocument

# --- Paragraph 1: "On Pilgrimage - June 1973" -------------------------------
# Heading1 -> Title (pandoc-style title block). Text stays the same.
$titlePara = $d.Paragraphs(1)
$titlePara.Style = "Title"

# --- Paragraph 2: "By Dorothy Day" -> "Dorothy Day" --------------------------
# Drop the old paragraph (and its direct-bold run) entirely, then grow a brand
# new paragraph right after the title text (but before the paragraph mark) so
# the new run inherits clean (non-bold) formatting instead of the stale
# <w:b/> the old "By Dorothy Day" run carried.
$byLinePara = $d.Paragraphs(2)
$byLinePara.Range.Delete()

$titleRange = $d.Paragraphs(1).Range
$splitPoint = $d.Range($titleRange.End - 1, $titleRange.End - 1)
$splitPoint.InsertParagraphAfter()

$authorsPara = $d.Paragraphs(2)
$authorsPara.Style = "Authors"
$authorsPara.Range.Text = "Dorothy Day"

Write-Output ("P1 text=[" + $d.Paragraphs(1).Range.Text + "] style=" + $d.Paragraphs(1).Style.NameLocal)
Write-Output ("P2 text=[" + $d.Paragraphs(2).Range.Text + "] style=" + $d.Paragraphs(2).Style.NameLocal + " bold=" + $d.Paragraphs(2).Range.Bold)
